$wb = $excel.ActiveWorkbook

# ----- Sheet "Fallecido_Recuperado": append 2 weekly summary rows -----
$ws1 = $wb.Worksheets.Item("Fallecido_Recuperado")

# Copy the date-format style from the last existing row (A34) onto the new date cells
$ws1.Range("A34").Copy() | Out-Null
$ws1.Range("A35:A36").PasteSpecial(-4122) | Out-Null

$ws1.Cells.Item(35,1).Value = 44149
$ws1.Cells.Item(35,2).Value = 135157
$ws1.Cells.Item(35,3).Value = 2293
$ws1.Cells.Item(35,4).Value = 110871
$ws1.Cells.Item(36,1).Value = 44156
$ws1.Cells.Item(36,2).Value = 138410
$ws1.Cells.Item(36,3).Value = 2310
$ws1.Cells.Item(36,4).Value = 112552

# ----- Sheet "Provincias_Semanal": append 2 full weeks x 32 provinces = 64 rows -----
$ws2 = $wb.Worksheets.Item("Provincias_Semanal")

# Copy the date-format style from the last existing row (A1057) onto all new date cells
$ws2.Range("A1057").Copy() | Out-Null
$ws2.Range("A1058:A1121").PasteSpecial(-4122) | Out-Null

$ws2.Cells.Item(1058,1).Value = 44149
$ws2.Cells.Item(1058,2).Value = "Distrito Nacional"
$ws2.Cells.Item(1058,3).Value = 3383.86
$ws2.Cells.Item(1058,4).Value = 383
$ws2.Cells.Item(1059,1).Value = 44149
$ws2.Cells.Item(1059,2).Value = "Azua"
$ws2.Cells.Item(1059,3).Value = 927.31
$ws2.Cells.Item(1059,4).Value = 24
$ws2.Cells.Item(1060,1).Value = 44149
$ws2.Cells.Item(1060,2).Value = "Baoruco"
$ws2.Cells.Item(1060,3).Value = 977.31
$ws2.Cells.Item(1060,4).Value = 8
$ws2.Cells.Item(1061,1).Value = 44149
$ws2.Cells.Item(1061,2).Value = "Barahona"
$ws2.Cells.Item(1061,3).Value = 1096.33
$ws2.Cells.Item(1061,4).Value = 22
$ws2.Cells.Item(1062,1).Value = 44149
$ws2.Cells.Item(1062,2).Value = "Dajabon"
$ws2.Cells.Item(1062,3).Value = 606.32000000000005
$ws2.Cells.Item(1062,4).Value = 7
$ws2.Cells.Item(1063,1).Value = 44149
$ws2.Cells.Item(1063,2).Value = "Duarte"
$ws2.Cells.Item(1063,3).Value = 1192.6500000000001
$ws2.Cells.Item(1063,4).Value = 116
$ws2.Cells.Item(1064,1).Value = 44149
$ws2.Cells.Item(1064,2).Value = "Elias Pina"
$ws2.Cells.Item(1064,3).Value = 372.23
$ws2.Cells.Item(1064,4).Value = 5
$ws2.Cells.Item(1065,1).Value = 44149
$ws2.Cells.Item(1065,2).Value = "El Seibo"
$ws2.Cells.Item(1065,3).Value = 584.98
$ws2.Cells.Item(1065,4).Value = 6
$ws2.Cells.Item(1066,1).Value = 44149
$ws2.Cells.Item(1066,2).Value = "Espaillat"
$ws2.Cells.Item(1066,3).Value = 1162.42
$ws2.Cells.Item(1066,4).Value = 89
$ws2.Cells.Item(1067,1).Value = 44149
$ws2.Cells.Item(1067,2).Value = "Independencia"
$ws2.Cells.Item(1067,3).Value = 1177.5999999999999
$ws2.Cells.Item(1067,4).Value = 9
$ws2.Cells.Item(1068,1).Value = 44149
$ws2.Cells.Item(1068,2).Value = "La Altagracia"
$ws2.Cells.Item(1068,3).Value = 1401.79
$ws2.Cells.Item(1068,4).Value = 41
$ws2.Cells.Item(1069,1).Value = 44149
$ws2.Cells.Item(1069,2).Value = "La Romana"
$ws2.Cells.Item(1069,3).Value = 1230.75
$ws2.Cells.Item(1069,4).Value = 49
$ws2.Cells.Item(1070,1).Value = 44149
$ws2.Cells.Item(1070,2).Value = "La Vega"
$ws2.Cells.Item(1070,3).Value = 1258.97
$ws2.Cells.Item(1070,4).Value = 115
$ws2.Cells.Item(1071,1).Value = 44149
$ws2.Cells.Item(1071,2).Value = "Maria Trinidad Sanchez"
$ws2.Cells.Item(1071,3).Value = 901.51
$ws2.Cells.Item(1071,4).Value = 12
$ws2.Cells.Item(1072,1).Value = 44149
$ws2.Cells.Item(1072,2).Value = "Monte Cristi"
$ws2.Cells.Item(1072,3).Value = 560.48
$ws2.Cells.Item(1072,4).Value = 15
$ws2.Cells.Item(1073,1).Value = 44149
$ws2.Cells.Item(1073,2).Value = "Pedernales"
$ws2.Cells.Item(1073,3).Value = 1611.57
$ws2.Cells.Item(1073,4).Value = 3
$ws2.Cells.Item(1074,1).Value = 44149
$ws2.Cells.Item(1074,2).Value = "Peravia"
$ws2.Cells.Item(1074,3).Value = 686.81
$ws2.Cells.Item(1074,4).Value = 45
$ws2.Cells.Item(1075,1).Value = 44149
$ws2.Cells.Item(1075,2).Value = "Puerto Plata"
$ws2.Cells.Item(1075,3).Value = 1084.57
$ws2.Cells.Item(1075,4).Value = 133
$ws2.Cells.Item(1076,1).Value = 44149
$ws2.Cells.Item(1076,2).Value = "Hermanas Mirabal"
$ws2.Cells.Item(1076,3).Value = 1071.0999999999999
$ws2.Cells.Item(1076,4).Value = 22
$ws2.Cells.Item(1077,1).Value = 44149
$ws2.Cells.Item(1077,2).Value = "Samana"
$ws2.Cells.Item(1077,3).Value = 487.7
$ws2.Cells.Item(1077,4).Value = 3
$ws2.Cells.Item(1078,1).Value = 44149
$ws2.Cells.Item(1078,2).Value = "San Cristobal"
$ws2.Cells.Item(1078,3).Value = 767.14
$ws2.Cells.Item(1078,4).Value = 117
$ws2.Cells.Item(1079,1).Value = 44149
$ws2.Cells.Item(1079,2).Value = "San Juan"
$ws2.Cells.Item(1079,3).Value = 1041.33
$ws2.Cells.Item(1079,4).Value = 43
$ws2.Cells.Item(1080,1).Value = 44149
$ws2.Cells.Item(1080,2).Value = "San Pedro de Macoris"
$ws2.Cells.Item(1080,3).Value = 584.33000000000004
$ws2.Cells.Item(1080,4).Value = 48
$ws2.Cells.Item(1081,1).Value = 44149
$ws2.Cells.Item(1081,2).Value = "Sanchez Ramirez"
$ws2.Cells.Item(1081,3).Value = 1342.3
$ws2.Cells.Item(1081,4).Value = 18
$ws2.Cells.Item(1082,1).Value = 44149
$ws2.Cells.Item(1082,2).Value = "Santiago"
$ws2.Cells.Item(1082,3).Value = 1408.67
$ws2.Cells.Item(1082,4).Value = 343
$ws2.Cells.Item(1083,1).Value = 44149
$ws2.Cells.Item(1083,2).Value = "Santiago Rodriguez"
$ws2.Cells.Item(1083,3).Value = 1169.98
$ws2.Cells.Item(1083,4).Value = 11
$ws2.Cells.Item(1084,1).Value = 44149
$ws2.Cells.Item(1084,2).Value = "Valverde"
$ws2.Cells.Item(1084,3).Value = 550.02
$ws2.Cells.Item(1084,4).Value = 31
$ws2.Cells.Item(1085,1).Value = 44149
$ws2.Cells.Item(1085,2).Value = "Monsenor Nouel"
$ws2.Cells.Item(1085,3).Value = 1086.2
$ws2.Cells.Item(1085,4).Value = 32
$ws2.Cells.Item(1086,1).Value = 44149
$ws2.Cells.Item(1086,2).Value = "Monte Plata"
$ws2.Cells.Item(1086,3).Value = 312.51
$ws2.Cells.Item(1086,4).Value = 28
$ws2.Cells.Item(1087,1).Value = 44149
$ws2.Cells.Item(1087,2).Value = "Hato Mayor"
$ws2.Cells.Item(1087,3).Value = 569.02
$ws2.Cells.Item(1087,4).Value = 13
$ws2.Cells.Item(1088,1).Value = 44149
$ws2.Cells.Item(1088,2).Value = "San Jose de Ocoa"
$ws2.Cells.Item(1088,3).Value = 1060.56
$ws2.Cells.Item(1088,4).Value = 13
$ws2.Cells.Item(1089,1).Value = 44149
$ws2.Cells.Item(1089,2).Value = "Santo Domingo"
$ws2.Cells.Item(1089,3).Value = 941.88
$ws2.Cells.Item(1089,4).Value = 489
$ws2.Cells.Item(1090,1).Value = 44156
$ws2.Cells.Item(1090,2).Value = "Distrito Nacional"
$ws2.Cells.Item(1090,3).Value = 3485.67
$ws2.Cells.Item(1090,4).Value = 386
$ws2.Cells.Item(1091,1).Value = 44156
$ws2.Cells.Item(1091,2).Value = "Azua"
$ws2.Cells.Item(1091,3).Value = 929.11
$ws2.Cells.Item(1091,4).Value = 24
$ws2.Cells.Item(1092,1).Value = 44156
$ws2.Cells.Item(1092,2).Value = "Baoruco"
$ws2.Cells.Item(1092,3).Value = 979.16
$ws2.Cells.Item(1092,4).Value = 8
$ws2.Cells.Item(1093,1).Value = 44156
$ws2.Cells.Item(1093,2).Value = "Barahona"
$ws2.Cells.Item(1093,3).Value = 1097.3800000000001
$ws2.Cells.Item(1093,4).Value = 22
$ws2.Cells.Item(1094,1).Value = 44156
$ws2.Cells.Item(1094,2).Value = "Dajabon"
$ws2.Cells.Item(1094,3).Value = 655.96
$ws2.Cells.Item(1094,4).Value = 7
$ws2.Cells.Item(1095,1).Value = 44156
$ws2.Cells.Item(1095,2).Value = "Duarte"
$ws2.Cells.Item(1095,3).Value = 1257.22
$ws2.Cells.Item(1095,4).Value = 116
$ws2.Cells.Item(1096,1).Value = 44156
$ws2.Cells.Item(1096,2).Value = "Elias Pina"
$ws2.Cells.Item(1096,3).Value = 375.38
$ws2.Cells.Item(1096,4).Value = 5
$ws2.Cells.Item(1097,1).Value = 44156
$ws2.Cells.Item(1097,2).Value = "El Seibo"
$ws2.Cells.Item(1097,3).Value = 586.04999999999995
$ws2.Cells.Item(1097,4).Value = 6
$ws2.Cells.Item(1098,1).Value = 44156
$ws2.Cells.Item(1098,2).Value = "Espaillat"
$ws2.Cells.Item(1098,3).Value = 1183.26
$ws2.Cells.Item(1098,4).Value = 90
$ws2.Cells.Item(1099,1).Value = 44156
$ws2.Cells.Item(1099,2).Value = "Independencia"
$ws2.Cells.Item(1099,3).Value = 1223.81
$ws2.Cells.Item(1099,4).Value = 9
$ws2.Cells.Item(1100,1).Value = 44156
$ws2.Cells.Item(1100,2).Value = "La Altagracia"
$ws2.Cells.Item(1100,3).Value = 1429.8
$ws2.Cells.Item(1100,4).Value = 42
$ws2.Cells.Item(1101,1).Value = 44156
$ws2.Cells.Item(1101,2).Value = "La Romana"
$ws2.Cells.Item(1101,3).Value = 1245.06
$ws2.Cells.Item(1101,4).Value = 49
$ws2.Cells.Item(1102,1).Value = 44156
$ws2.Cells.Item(1102,2).Value = "La Vega"
$ws2.Cells.Item(1102,3).Value = 1307.1099999999999
$ws2.Cells.Item(1102,4).Value = 117
$ws2.Cells.Item(1103,1).Value = 44156
$ws2.Cells.Item(1103,2).Value = "Maria Trinidad Sanchez"
$ws2.Cells.Item(1103,3).Value = 930.57
$ws2.Cells.Item(1103,4).Value = 12
$ws2.Cells.Item(1104,1).Value = 44156
$ws2.Cells.Item(1104,2).Value = "Monte Cristi"
$ws2.Cells.Item(1104,3).Value = 561.33000000000004
$ws2.Cells.Item(1104,4).Value = 15
$ws2.Cells.Item(1105,1).Value = 44156
$ws2.Cells.Item(1105,2).Value = "Pedernales"
$ws2.Cells.Item(1105,3).Value = 1611.57
$ws2.Cells.Item(1105,4).Value = 3
$ws2.Cells.Item(1106,1).Value = 44156
$ws2.Cells.Item(1106,2).Value = "Peravia"
$ws2.Cells.Item(1106,3).Value = 689.85
$ws2.Cells.Item(1106,4).Value = 45
$ws2.Cells.Item(1107,1).Value = 44156
$ws2.Cells.Item(1107,2).Value = "Puerto Plata"
$ws2.Cells.Item(1107,3).Value = 1139.18
$ws2.Cells.Item(1107,4).Value = 133
$ws2.Cells.Item(1108,1).Value = 44156
$ws2.Cells.Item(1108,2).Value = "Hermanas Mirabal"
$ws2.Cells.Item(1108,3).Value = 1078.7
$ws2.Cells.Item(1108,4).Value = 23
$ws2.Cells.Item(1109,1).Value = 44156
$ws2.Cells.Item(1109,2).Value = "Samana"
$ws2.Cells.Item(1109,3).Value = 499.29
$ws2.Cells.Item(1109,4).Value = 3
$ws2.Cells.Item(1110,1).Value = 44156
$ws2.Cells.Item(1110,2).Value = "San Cristobal"
$ws2.Cells.Item(1110,3).Value = 778.28
$ws2.Cells.Item(1110,4).Value = 117
$ws2.Cells.Item(1111,1).Value = 44156
$ws2.Cells.Item(1111,2).Value = "San Juan"
$ws2.Cells.Item(1111,3).Value = 1049.9000000000001
$ws2.Cells.Item(1111,4).Value = 44
$ws2.Cells.Item(1112,1).Value = 44156
$ws2.Cells.Item(1112,2).Value = "San Pedro de Macoris"
$ws2.Cells.Item(1112,3).Value = 604.99
$ws2.Cells.Item(1112,4).Value = 49
$ws2.Cells.Item(1113,1).Value = 44156
$ws2.Cells.Item(1113,2).Value = "Sanchez Ramirez"
$ws2.Cells.Item(1113,3).Value = 1396.91
$ws2.Cells.Item(1113,4).Value = 18
$ws2.Cells.Item(1114,1).Value = 44156
$ws2.Cells.Item(1114,2).Value = "Santiago"
$ws2.Cells.Item(1114,3).Value = 1430.68
$ws2.Cells.Item(1114,4).Value = 344
$ws2.Cells.Item(1115,1).Value = 44156
$ws2.Cells.Item(1115,2).Value = "Santiago Rodriguez"
$ws2.Cells.Item(1115,3).Value = 1185.69
$ws2.Cells.Item(1115,4).Value = 11
$ws2.Cells.Item(1116,1).Value = 44156
$ws2.Cells.Item(1116,2).Value = "Valverde"
$ws2.Cells.Item(1116,3).Value = 552.85
$ws2.Cells.Item(1116,4).Value = 31
$ws2.Cells.Item(1117,1).Value = 44156
$ws2.Cells.Item(1117,2).Value = "Monsenor Nouel"
$ws2.Cells.Item(1117,3).Value = 1152.18
$ws2.Cells.Item(1117,4).Value = 32
$ws2.Cells.Item(1118,1).Value = 44156
$ws2.Cells.Item(1118,2).Value = "Monte Plata"
$ws2.Cells.Item(1118,3).Value = 316.18
$ws2.Cells.Item(1118,4).Value = 28
$ws2.Cells.Item(1119,1).Value = 44156
$ws2.Cells.Item(1119,2).Value = "Hato Mayor"
$ws2.Cells.Item(1119,3).Value = 571.35
$ws2.Cells.Item(1119,4).Value = 13
$ws2.Cells.Item(1120,1).Value = 44156
$ws2.Cells.Item(1120,2).Value = "San Jose de Ocoa"
$ws2.Cells.Item(1120,3).Value = 1062.4000000000001
$ws2.Cells.Item(1120,4).Value = 13
$ws2.Cells.Item(1121,1).Value = 44156
$ws2.Cells.Item(1121,2).Value = "Santo Domingo"
$ws2.Cells.Item(1121,3).Value = 964.62
$ws2.Cells.Item(1121,4).Value = 495

# ----- Restore view state (selection / scroll / active sheet) to match the edited workbook -----
$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$ws1.Range("E36").Select() | Out-Null

$ws2.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 1096
$ws2.Range("D1122").Select() | Out-Null

Write-Host "Applied: 2 new weekly rows (Fallecido_Recuperado) + 64 new provincial rows (Provincias_Semanal)"
